# Update "想去人数" (want-to-go count) values in both the "展览" (Exhibition)
# and "全部类型" (All types) sheets, which carry duplicated data rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 31
    $ws.Range("F3").Value = 64
}
